$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 25: add the end time (D25), which lets E25/F25 (and all the
# downstream weekly-summary formulas in row 4/5/7) recompute automatically.
# Copy D24's number-format (time) down first so D25 doesn't inherit a
# generic style when we stuff a Value into it.
$ws.Range("D24").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 0.96527777777777779

# --- Row 25: daily log text. Set H25 (next steps) before G25 (what got
# done) so the two new shared strings land in the same order the diff
# expects (76 = "I think the next thing...", 77 = "Refactored the state...").
$ws.Range("H24").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("H25").Value = "I think the next thing on the quest to refactor is addressing all clouds/thought bubbles 💭"

$ws.Range("G24").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G25").Value = "Refactored the state of the application to be less of a massive headache. Implemented custom hooks to grab select groups of values from the context. I also moved fragile code to a context."

# Row grows taller to fit the newly-wrapped text (matches the other
# journal rows, which are all auto-sized to their text).
$ws.Rows.Item(25).RowHeight = 90

# Recalculate everything so the SUMIFS/array-formula summary cells pick
# up the new row-25 duration.
$excel.CalculateFullRebuild()

# --- View state: scroll down one row and move the selection to G26,
# matching where the author's cursor ended up after typing the entry.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 19
$ws.Range("G26").Select()
